$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Position" column (D) currently hold the text "stock" and
# must instead hold the numeric value 1 (destination was recorded before
# the item actually reached its position).
$rows = @(3,5,6,9,10,11,12,13,14,17,18,19,20,22,23,24,25,26,27,28,29,30,31,32,33)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Move the active selection back to the top of the Position column (D1)
# instead of where it was left (D9).
$ws.Range("D1").Select()
